# Auto-generated edit script applying cached market-price snapshot updates
# from the scheduled runner, covering the ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 5000
$ws.Range("J7").Value = 5000
$ws.Range("L7").Value = 5000
$ws.Range("N7").Value = -5224
$ws.Range("H14").Value = 5000
$ws.Range("J14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("N14").Value = -5382
$ws.Range("H18").Value = 3416.3333
$ws.Range("I18").Value = 4499.6665
$ws.Range("J18").Value = 2333
$ws.Range("K18").Value = 4499.6665
$ws.Range("L18").Value = 2333
$ws.Range("M18").Value = -4215.6665
$ws.Range("N18").Value = -2901
$ws.Range("H62").Value = 2458
$ws.Range("J62").Value = 1903
$ws.Range("L62").Value = 1903
$ws.Range("N62").Value = -3151
$ws.Range("H65").Value = 2458
$ws.Range("J65").Value = 1903
$ws.Range("L65").Value = 9515
$ws.Range("N65").Value = -15755

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6772.9414
$ws.Range("I2").Value = 908
$ws.Range("J2").Value = 25834
$ws.Range("K2").Value = 908
$ws.Range("L2").Value = 25834
$ws.Range("M2").Value = -795
$ws.Range("N2").Value = -26060
$ws.Range("H32").Value = 4545.4194
$ws.Range("I32").Value = 5145.3076
$ws.Range("J32").Value = 1426
$ws.Range("K32").Value = 5145.3076
$ws.Range("L32").Value = 1426
$ws.Range("M32").Value = -4858.3076
$ws.Range("N32").Value = -2000
$ws.Range("H61").Value = 1912.2727
$ws.Range("I61").Value = 1565.125
$ws.Range("J61").Value = 2838
$ws.Range("K61").Value = 1565.125
$ws.Range("L61").Value = 2838
$ws.Range("M61").Value = -1353.125
$ws.Range("N61").Value = -3262
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H116").Value = 6772.9414
$ws.Range("I116").Value = 908
$ws.Range("J116").Value = 25834
$ws.Range("K116").Value = 908
$ws.Range("L116").Value = 25834
$ws.Range("M116").Value = 1386
$ws.Range("N116").Value = -30422
$ws.Range("H132").Value = 2423.125
$ws.Range("I132").Value = 2207.5264
$ws.Range("K132").Value = 6622.5792
$ws.Range("M132").Value = -4092.5792
$ws.Range("H136").Value = 1912.2727
$ws.Range("I136").Value = 1565.125
$ws.Range("J136").Value = 2838
$ws.Range("K136").Value = 4695.375
$ws.Range("L136").Value = 8514
$ws.Range("M136").Value = -2145.375
$ws.Range("N136").Value = -13614

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6772.9414
$ws.Range("I3").Value = 908
$ws.Range("J3").Value = 25834
$ws.Range("K3").Value = 908
$ws.Range("L3").Value = 25834
$ws.Range("M3").Value = -794
$ws.Range("N3").Value = -26062
$ws.Range("H5").Value = 505
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 505
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 505
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -731
$ws.Range("H22").Value = 680
$ws.Range("I22").Value = 680
$ws.Range("K22").Value = 680
$ws.Range("M22").Value = -507
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H105").Value = 200002460
$ws.Range("I105").Value = 333335680
$ws.Range("K105").Value = 333335680
$ws.Range("M105").Value = -333333933
$ws.Range("H107").Value = 1952.4736
$ws.Range("I107").Value = 1435.6666
$ws.Range("J107").Value = 2417.6
$ws.Range("K107").Value = 1435.6666
$ws.Range("L107").Value = 2417.6
$ws.Range("M107").Value = 484.3334
$ws.Range("N107").Value = -6257.6

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 319.7143
$ws.Range("I7").Value = 409.2
$ws.Range("K7").Value = 409.2
$ws.Range("M7").Value = -296.2
$ws.Range("H22").Value = 549.8333
$ws.Range("I22").Value = 560
$ws.Range("K22").Value = 560
$ws.Range("M22").Value = -210
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H31").Value = 1907.4615
$ws.Range("I31").Value = 920.7368
$ws.Range("K31").Value = 920.7368
$ws.Range("M31").Value = -625.7368
$ws.Range("H34").Value = 1907.4615
$ws.Range("I34").Value = 920.7368
$ws.Range("K34").Value = 920.7368
$ws.Range("M34").Value = -718.7368
$ws.Range("H58").Value = 1254.7142
$ws.Range("I58").Value = 1313.8334
$ws.Range("J58").Value = 900
$ws.Range("K58").Value = 1313.8334
$ws.Range("L58").Value = 900
$ws.Range("M58").Value = -1110.8334
$ws.Range("N58").Value = -1306
$ws.Range("H122").Value = 891.6667
$ws.Range("I122").Value = 712.5
$ws.Range("K122").Value = 2137.5
$ws.Range("M122").Value = 312.5
$ws.Range("H136").Value = 1254.7142
$ws.Range("I136").Value = 1313.8334
$ws.Range("J136").Value = 900
$ws.Range("K136").Value = 3941.5002
$ws.Range("L136").Value = 2700
$ws.Range("M136").Value = -1391.5002
$ws.Range("N136").Value = -7800

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2286.7144
$ws.Range("I69").Value = 600
$ws.Range("J69").Value = 2416.4614
$ws.Range("K69").Value = 1800
$ws.Range("L69").Value = 7249.3842
$ws.Range("M69").Value = -989
$ws.Range("N69").Value = -8871.3842
$ws.Range("H72").Value = 2286.7144
$ws.Range("I72").Value = 600
$ws.Range("J72").Value = 2416.4614
$ws.Range("K72").Value = 5400
$ws.Range("L72").Value = 21748.1526
$ws.Range("M72").Value = -1344
$ws.Range("N72").Value = -29860.1526
$ws.Range("H97").Value = 613.8570999999999
$ws.Range("I97").Value = 599.3333
$ws.Range("J97").Value = 624.75
$ws.Range("K97").Value = 1797.9999
$ws.Range("L97").Value = 1874.25
$ws.Range("M97").Value = -1301.9999
$ws.Range("N97").Value = -2866.25

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 21434614
$ws.Range("I70").Value = 20838518
$ws.Range("K70").Value = 20838518
$ws.Range("M70").Value = -20838248
$ws.Range("H73").Value = 21434614
$ws.Range("I73").Value = 20838518
$ws.Range("K73").Value = 20838518
$ws.Range("M73").Value = -20837582
$ws.Range("H132").Value = 2582.3157
$ws.Range("I132").Value = 2191.0667
$ws.Range("J132").Value = 4049.5
$ws.Range("K132").Value = 6573.2001
$ws.Range("L132").Value = 12148.5
$ws.Range("M132").Value = -4043.2001
$ws.Range("N132").Value = -17208.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2176
$ws.Range("I7").Value = 2150
$ws.Range("J7").Value = 2186.4
$ws.Range("K7").Value = 2150
$ws.Range("L7").Value = 2186.4
$ws.Range("M7").Value = -2038
$ws.Range("N7").Value = -2410.4
$ws.Range("H46").Value = 2174.8572
$ws.Range("I46").Value = 1944.4
$ws.Range("K46").Value = 1944.4
$ws.Range("M46").Value = -1756.4
$ws.Range("H96").Value = 30197
$ws.Range("J96").Value = 30197
$ws.Range("L96").Value = 30197
$ws.Range("N96").Value = -35689
$ws.Range("H126").Value = 2176
$ws.Range("I126").Value = 2150
$ws.Range("J126").Value = 2186.4
$ws.Range("K126").Value = 6450
$ws.Range("L126").Value = 6559.200000000001
$ws.Range("M126").Value = -3980
$ws.Range("N126").Value = -11499.2
$ws.Range("H136").Value = 6805.8423
$ws.Range("I136").Value = 8357.929
$ws.Range("J136").Value = 2460
$ws.Range("K136").Value = 25073.787
$ws.Range("L136").Value = 7380
$ws.Range("M136").Value = -22523.787
$ws.Range("N136").Value = -12480

